$p = $ppt.ActivePresentation

# --- Slide 16 ("Prim's MST Algorithm"): left-align the
#     "Can this lead to an infeasible solution?" paragraph ---
$slide16 = $p.Slides.Item(16)
$shape16 = $slide16.Shapes.Item(3)
$tr16 = $shape16.TextFrame.TextRange
$para16 = $tr16.Paragraphs(7, 1)
$para16.ParagraphFormat.Alignment = 1

# --- Slide 17 ("Tracking Edges for Prim's MST"): fix typo
#     "Candidates edges" -> "Candidate edges" ---
$slide17 = $p.Slides.Item(17)
$shape17 = $slide17.Shapes.Item(3)
$tr17 = $shape17.TextFrame.TextRange
$para17 = $tr17.Paragraphs(1, 1)
$run17 = $tr17.Characters($para17.Start, $para17.Length)
$run17.Text = "Candidate edges:  edge from a tree-node to a non-tree node"
